# "updated NTT Address on the first page"
#
# The title-page address block had "Cluj Napoca" in its own paragraph,
# followed by a second paragraph "City, State 400158" (two runs:
# "City, State " and "400158"). The edit merges both paragraphs into a
# single paragraph containing one run with the text
# "Cluj Napoca City, 400158" (dropping the word "State" and joining the
# two lines), reusing the formatting (bold, small caps, 28pt) of the
# original "Cluj Napoca" run.

$d = $word.ActiveDocument

# Find the paragraph break (^p) between "Cluj Napoca" and
# "City, State 400158" and replace the whole span with the new,
# single-line address text. This collapses the two paragraphs into one
# and leaves a single run behind, carrying the formatting of the first
# matched run.
$found = $d.Content.Find.Execute( `
    "Cluj Napoca^pCity, State 400158", `
    $true, `
    $false, `
    $false, `
    $false, `
    $false, `
    $true, `
    1, `
    $false, `
    "Cluj Napoca City, 400158", `
    2)

if (-not $found) {
    throw "Could not find the 'Cluj Napoca' / 'City, State 400158' address paragraphs to merge."
}
